$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(3).Insert()
$ws.Range("C1").Value = "Rating"
$ws.Range("G1").Value = "Website"

$ws.Hyperlinks.Add($ws.Range("G2"), "http://www.indianapoliszoo.com/")

Write-Host "done"
